# "Generate Report for Handback" - mark the two localization entries as
# handed back (in sync with en-US) on both the zh-cn and de-de sheets,
# filling in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns and linking the new target-file cell
# back to the source markdown file, just like column A already does.

$wb = $excel.ActiveWorkbook

$mdUrl00ea = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/73f21880514873b461e36b9e2e31551e9af5ec6f/e2e/00ea13d2-1419-4f8c-9d94-be771667e003.md"
$mdUrl51bc = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/73f21880514873b461e36b9e2e31551e9af5ec6f/e2e/51bcdee1-b9d0-4a67-a8b4-136896f59d72.md"

$md00ea = "00ea13d2-1419-4f8c-9d94-be771667e003.md"
$md51bc = "51bcdee1-b9d0-4a67-a8b4-136896f59d72.md"

$handedBackStatus = "Handed back: in sync with en-US"

# hyperlink font color used elsewhere in the workbook (RGB FF6495ED as an OLE BGR long)
$hyperlinkColor = 15570276

function Set-HandbackRow {
    # positional params: worksheet, row, md display name, md url, xlf name, handback datetime
    param($ws, $row, $mdName, $mdUrl, $xlfName, $handbackDateTime)

    # Status -> Handed back
    $ws.Cells.Item($row, 3).Value = $handedBackStatus

    # Latest Target File (col I) -- becomes a hyperlink to the source .md,
    # same as column A's hyperlink.
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 9), $mdUrl, "", "", $mdName)
    $ws.Cells.Item($row, 9).Font.Underline = 2
    $ws.Cells.Item($row, 9).Font.Color = $hyperlinkColor

    # Latest Handback File (col J) -- the xlf that was handed back.
    $ws.Cells.Item($row, 10).Value = $xlfName

    # Latest Handback DateTime (col K)
    $ws.Cells.Item($row, 11).Value = $handbackDateTime
}

# The "Overview" sheet has its own Status columns (E = zh-cn, F = de-de)
# that show the very same text, so flip those too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 5).Value = $handedBackStatus
$wsOverview.Cells.Item(2, 6).Value = $handedBackStatus
$wsOverview.Cells.Item(3, 5).Value = $handedBackStatus
$wsOverview.Cells.Item(3, 6).Value = $handedBackStatus

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "zh-cn") {
        $handbackTime = "2016-08-16 23:00:32"
        $xlf00ea = "00ea13d2-1419-4f8c-9d94-be771667e003.064b06a2cb5c5f99c40fbdc681e8a627cbce8380.zh-cn.xlf"
        $xlf51bc = "51bcdee1-b9d0-4a67-a8b4-136896f59d72.989e6333fa360a38a15ecab03460da3f55733caf.zh-cn.xlf"
    } else {
        $handbackTime = "2016-08-16 23:00:40"
        $xlf00ea = "00ea13d2-1419-4f8c-9d94-be771667e003.064b06a2cb5c5f99c40fbdc681e8a627cbce8380.de-de.xlf"
        $xlf51bc = "51bcdee1-b9d0-4a67-a8b4-136896f59d72.989e6333fa360a38a15ecab03460da3f55733caf.de-de.xlf"
    }

    Set-HandbackRow $ws 2 $md00ea $mdUrl00ea $xlf00ea $handbackTime
    Set-HandbackRow $ws 3 $md51bc $mdUrl51bc $xlf51bc $handbackTime

    # Widen the Status / Latest Target File / Latest Handback File columns
    # so the longer text/hyperlinks fit.
    $ws.Columns.Item(3).ColumnWidth = 175 / 6
    $ws.Columns.Item(9).ColumnWidth = 235 / 6
    $ws.Columns.Item(10).ColumnWidth = 235 / 6
}
